$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Observed"
$ws.Name = "Observed"

# The new shared strings must be created in this exact order so that the
# underlying shared-strings table ends up with:
#   14 = Sim2, 15 = Sim3, 16 = Sim5, 17 = Simulation1
# (matching the target workbook's sharedStrings.xml ordering)
$ws.Range("A3").Value = "Sim2"
$ws.Range("A4").Value = "Sim3"
$ws.Range("A7").Value = "Sim5"
$ws.Range("A2").Value = "Simulation1"

# Re-use the already created shared strings for the remaining rows that
# reference the same simulation names.
$ws.Range("A5").Value = "Simulation1"
$ws.Range("A6").Value = "Simulation1"
$ws.Range("A8").Value = "Sim2"

# New numeric values
$ws.Range("F2").Value = 10

$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 20

$ws.Range("B8").Value = "2/01/2000"
$ws.Range("G8").Value = 20

# Update the selected cell shown in the saved worksheet view
[void]$ws.Range("E20").Select()
